# DOMA-4173: add new "order" column into the payments export template.
#
# A new column is inserted before column G ("status"), pushing the existing
# "status"/"amount" columns one place to the right (to H/I), and the new
# column G is filled with the "order" placeholders for the header row and
# the two sample/body (item) rows, matching the formatting of the column
# immediately to its left ("transaction").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at G; existing G ("status") and H ("amount") columns
# shift one place to the right (to H and I).
$ws.Columns("G:G").Insert()

# Copy formatting (styles/borders/fills) from column F into the freshly
# inserted column G so it matches the look of the surrounding header/body
# cells.
$ws.Range("F1:F10").Copy()
$ws.Range("G1:G10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Widen the new "order" column.
$ws.Columns("G:G").ColumnWidth = 25.1

# New "order" placeholders for the header row and the two item rows.
$ws.Range("G1").Value = "{d.i18n.order}"
$ws.Range("G2").Value = "{d.objs[I].order}"
$ws.Range("G3").Value = "{d.objs[I+1].order}"
